$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post "「スマイル」ابتسام..." that occupied row 563 was removed from the
# source data; every following row shifts up by one. Deleting the entire
# row 563 reproduces this (Excel shifts rows 564:632 up to 563:631).
$ws.Rows(563).Delete()
